$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.583.99"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "'3.047.68"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'534.40"
$ws.Range("E5").Value = "  -3.94%  "
$ws.Range("D6").Value = "'132.55"
$ws.Range("E6").Value = "  -4.43%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'3.036.61"
$ws.Range("E8").Value = "  -1.83%  "
$ws.Range("D9").Value = "'0.495"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").Value = "'0.153"
$ws.Range("E10").Value = "  +1.15%  "
$ws.Range("D11").Value = "'6.16"
$ws.Range("E11").Value = "  -8.67%  "
$ws.Range("D12").Value = "'0.451"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("D13").Value = "'0.0000223"
$ws.Range("E13").Value = "  +3.46%  "
$ws.Range("D14").Value = "'34.12"
$ws.Range("E14").Value = "  -3.80%  "
$ws.Range("D15").Value = "'3.555.22"
$ws.Range("D16").Value = "'62.646.39"
$ws.Range("E16").Value = "  -1.65%  "
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").Value = "'3.057.23"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").Value = "'6.58"
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("D20").Value = "'481.06"
$ws.Range("E20").Value = "  -4.68%  "
$ws.Range("D21").Value = "'13.20"
$ws.Range("E21").Value = "  -3.71%  "
$ws.Range("D22").Value = "'0.691"
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("D23").Value = "'7.09"
$ws.Range("E23").Value = "  -2.11%  "
$ws.Range("D24").Value = "'79.00"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").Value = "'12.09"
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "'2.69"
$ws.Range("E27").Value = "  -2.90%  "
$ws.Range("D28").Value = "'8.05"
$ws.Range("E28").Value = "  -3.95%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'25.79"
$ws.Range("E30").Value = "  -1.65%  "
$ws.Range("D31").Value = "'1.85"
$ws.Range("E31").Value = "  -9.47%  "
$ws.Range("D32").Value = "'1.10"
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("D33").Value = "'2.36"
$ws.Range("E33").Value = "  -7.12%  "
$ws.Range("D34").Value = "'56.32"
$ws.Range("E34").Value = "  +2.39%  "
$ws.Range("D35").Value = "'5.37"
$ws.Range("E35").Value = "  +3.38%  "
$ws.Range("D36").Value = "'5.91"
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("D37").Value = "'474.90"
$ws.Range("E37").Value = "  -12.26%  "
$ws.Range("D38").Value = "'0.0394"
$ws.Range("E38").Value = "  -4.76%  "
$ws.Range("D39").Value = "'3.075.75"
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").Value = "'0.0792"
$ws.Range("E40").Value = "  -1.03%  "
$ws.Range("D41").Value = "'0.115"
$ws.Range("E41").Value = "  -3.07%  "
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").Value = "'8.05"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.63"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("D44").Value = "'0.252"
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("D46").Value = "'0.0₃0544"
$ws.Range("E46").Value = "  +8.67%  "
$ws.Range("D47").Value = "'2.02"
$ws.Range("E47").Value = "  -3.78%  "
$ws.Range("D48").Value = "'120.42"
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("D49").Value = "'24.55"
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("D50").Value = "'0.107"
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("D51").Value = "'2.30"
$ws.Range("E51").Value = "  +2.17%  "
